# "Generate Report for Handback"
#
# - Status text "Ready for handoff" -> "Handback transform failed"
#   (every cell that shows that status, across all three sheets)
# - Error Detail column (P) widened from ~13.75 chars to 40 chars on the
#   zh-cn and de-de report sheets
# - Error Detail cell (P3) on each locale sheet now carries a
#   "Handback file name ... is different with handoff file name ..."
#   diagnostic message instead of being blank

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text update -------------------------------------------------
# All cells that currently read "Ready for handoff" need to move to the
# new status in lockstep so they keep sharing one string.
$newStatus = "Handback transform failed"

$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Error Detail column width ------------------------------------------
# Raw OOXML column width of 40 characters corresponds to a COM
# ColumnWidth of (40 - 0.8333333333333334); the runtime adds back the
# standard 5px/MDW padding when it serialises the sheet.
$errorDetailWidth = 40 - 0.8333333333333334
$wsZhCn.Columns.Item(16).ColumnWidth = $errorDetailWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $errorDetailWidth

# --- Error Detail messages ------------------------------------------------
$wsZhCn.Range("P3").Value = "Handback file name: 3w24v3pw.b2r is different with handoff file name: 042e2f9a-1f81-4830-9f4f-781bf1b3d165.a32cf37dfd6c76406f1b913becaa3bdcd4d5b0f3.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: 3w24v3pw.b2r is different with handoff file name: 042e2f9a-1f81-4830-9f4f-781bf1b3d165.a32cf37dfd6c76406f1b913becaa3bdcd4d5b0f3.de-de."
